# Update the "persistent setting" column header on the options sheet:
# shorten the text and reduce row 1's height to fit the new (shorter)
# heading, then move the active selection to C7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G1 held "Setting is persistent (stored in database file)" - shorten it.
$ws.Range("G1").Value = "Persistent (stored in database file)"

# Row 1 can now be shorter since the header text wraps less.
$ws.Rows.Item(1).RowHeight = 59.25

# Move the selection/active cell to C7.
$ws.Range("C7").Select()
